# The workbook was last saved with the "Formula" sheet (the 5th tab) active
# and the "raws" sheet scrolled down to row 31. The author switched back to
# the first tab ("raws") before saving, which:
#   - makes "raws" the active/selected sheet (tabSelected on its sheetView,
#     activeTab on the workbook's bookViews),
#   - clears the "Formula" sheet's tabSelected flag,
#   - resets the "raws" view back to the top-left (A1) instead of the
#     previously scrolled position.
# A full recalculation also happens on save (the model contains volatile
# RAND()/NORM.INV formulas that drive both the raw data table and the
# dependent charts), which is why cached values throughout the workbook and
# the chart data caches differ from the previous save.

$wb = $excel.ActiveWorkbook

$raws = $wb.Worksheets.Item("raws")
$raws.Activate()
